$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "PEÇAS" content row in the parts table: collapse the 4 runs
#    ("             " bold, "{part}" plain, " " bold, "{#parts}" plain)
#    into " {part} {#parts}" (plain) followed by "         " (9 spaces,
#    bold) -- matches the new layout from the diff.
# ------------------------------------------------------------------
$partsTable = $d.Tables.Item(4)
$partsCell = $partsTable.Rows.Item(2).Cells.Item(1)
$partsRange = $partsCell.Range
$partsStart = $partsRange.Start

$oldVisibleLen = 28  # "             {part} {#parts}".Length
$target = $d.Range($partsStart, $partsStart + $oldVisibleLen)
$newText = " {part} {#parts}         "
$target.Text = $newText

$newLen = $newText.Length
$wholeNew = $d.Range($partsStart, $partsStart + $newLen)
$wholeNew.Font.Bold = 0

$tailLen = 9
$tailRange = $d.Range($partsStart + $newLen - $tailLen, $partsStart + $newLen)
$tailRange.Font.Bold = 1

# ------------------------------------------------------------------
# 2) Shading: the "R$ {value}" cell goes from cccccc to b7b7b7
#    (matching the rest of that row).
# ------------------------------------------------------------------
$valueCell = $partsTable.Rows.Item(3).Cells.Item(2)
$valueCell.Shading.BackgroundPatternColor = 0xb7b7b7

# ------------------------------------------------------------------
# 3) Shading: the "{/parts}" footer row goes from 999999 to d9d9d9.
# ------------------------------------------------------------------
$footerCell = $partsTable.Rows.Item(4).Cells.Item(1)
$footerCell.Shading.BackgroundPatternColor = 0xd9d9d9

# ------------------------------------------------------------------
# 4) "VALOR TOTAL" value cell: "{totalValue}" becomes "R$  {totalValue}".
# ------------------------------------------------------------------
$totalTable = $d.Tables.Item(5)
$totalCell = $totalTable.Rows.Item(1).Cells.Item(2)
$totalStart = $totalCell.Range.Start
$firstChar = $d.Range($totalStart, $totalStart + 1)
$firstChar.InsertBefore("R$  ")
